$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "i-0ab7e42518e7a90f0"
$ws.Range("B6").Value = "t2.micro"
$ws.Range("D6").Value = "172.31.0.189"
$ws.Range("E6").Value = "2023-11-17 08:53:58+00:00"
$ws.Range("F6").Value = "launch-wizard-2"
